# feat: add 2022-Q4 data
#
# The workbook tracks BAYN (Bayer) fund-holding data per quarter. Each
# quarter gets its own worksheet (named e.g. "2022-Q3") plus a summary
# row on the "总计" (Total) worksheet. This script adds a new "2022-Q4"
# quarter: a new worksheet with that quarter's fund data, and a new
# top row on the summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet ("总计"): insert a new row 2 for 2022-Q4 and shift
#    the existing quarters down (renumbering the running index in A).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.18

for ($r = 3; $r -le 8; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet: clone the "2022-Q3" sheet (same column
#    layout/styling/fund name) right before it, rename the clone, and
#    overwrite the quarter-specific figures.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)

$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# D:G hold percentages/amounts stored as text (matching the rest of the
# workbook's convention), so force a text number-format before writing
# them or Excel would silently coerce "92.90" -> 92.9 (numeric).
$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "4.76"
$q4.Range("E2").Value = "92.90"
$q4.Range("F2").Value = "3.80"
$q4.Range("G2").Value = "0.1809"
$q4.Range("H2").Value = 8
